$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.746.73'
$ws.Range("E2").Value = '  +2.04%  '

$ws.Range("D3").Value = '1.865.38'
$ws.Range("E3").Value = '  +1.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.037'
$ws.Range("E4").Value = '  +2.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.41'
$ws.Range("E5").Value = '  +2.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.033'
$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4420'
$ws.Range("E7").Value = '  +2.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3796'
$ws.Range("E8").Value = '  +2.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07474'
$ws.Range("E9").Value = '  +2.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8842'
$ws.Range("E10").Value = '  +1.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.77'
$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("D12").Value = '1.887.03'
$ws.Range("E12").Value = '  -13.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.550'
$ws.Range("E13").Value = '  +1.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.761'
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07224'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.26'
$ws.Range("E16").Value = '  +3.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.039'
$ws.Range("E17").Value = '  +1.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009131'
$ws.Range("E18").Value = '  +1.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.033'
$ws.Range("E19").Value = '  +1.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.57'
$ws.Range("E20").Value = '  +1.88%  '

$ws.Range("D21").Value = '27.759.97'
$ws.Range("E21").Value = '  +1.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.314'
$ws.Range("E22").Value = '  +2.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.31'
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").Value = '2.101.11'
$ws.Range("E24").Value = '  -11.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.012'
$ws.Range("E25").Value = '  +6.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.62'
$ws.Range("E26").Value = '  +2.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.84'
$ws.Range("E27").Value = '  +2.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.994'
$ws.Range("E28").Value = '  +4.34%  '

$ws.Range("E29").Value = '  +1.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.97'
$ws.Range("E30").Value = '  +2.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09058'
$ws.Range("E31").Value = '  +0.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7791'
$ws.Range("E32").Value = '  +2.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.217'
$ws.Range("E33").Value = '  +1.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.019'
$ws.Range("E34").Value = '  +6.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.579'
$ws.Range("E35").Value = '  +2.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.035'
$ws.Range("E36").Value = '  +1.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01989'
$ws.Range("E38").Value = '  +2.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05341'
$ws.Range("E39").Value = '  +1.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.874'
$ws.Range("E40").Value = '  +3.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5206'
$ws.Range("E41").Value = '  +1.31%  '

$ws.Range("E42").Value = '  +2.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.866'
$ws.Range("E43").Value = '  +5.72%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.658'
$ws.Range("E44").Value = '  +3.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.35'
$ws.Range("E45").Value = '  +2.63%  '

$ws.Range("E46").Value = '  +2.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06614'
$ws.Range("E47").Value = '  +5.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.718'
$ws.Range("E48").Value = '  +3.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4717'
$ws.Range("E49").Value = '  +2.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.927'
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.75'
$ws.Range("E51").Value = '  +1.87%  '
